# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# The source rows for several fixtures were swapped (match data ended up on
# the wrong id/row). This fixes it by swapping the full contents of columns
# B..AC between the following row pairs, while leaving column A (the row's
# sequential id) untouched:
#   22 <-> 23
#   78 <-> 79
#   103 <-> 104
#   135 <-> 136

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, $rowA, $rowB, $colStart, $colEnd)
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Columns B (2) through AC (29)
Swap-RowData $ws 22 23 2 29
Swap-RowData $ws 78 79 2 29
Swap-RowData $ws 103 104 2 29
Swap-RowData $ws 135 136 2 29
